# edit.ps1 - apply "update information for dgl domain apps" changes
# - bump the cached datetimeFigureOut field text from 3/27/20 to 4/6/20
#   on the slide master, notes master and every slide layout
# - expand the two GCMC bullet lines on slide 3
# - relabel the u0/i0 node boxes on slide 4 to u1..u7 / i1..i7 (one of
#   each group stays "u0"/"i0")

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder ("3/27/20" -> "4/6/20") on master, notes master,
#    and all slide layouts.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.HasTextFrame) {
            $txt = $shp.TextFrame.TextRange.Text
            if ($txt -eq "3/27/20") {
                $shp.TextFrame.TextRange.Text = "4/6/20"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 3 - "GCMC" textbox: expand the two numbered steps.
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$gcmcBox = $slide3.Shapes.Item(1)
$gcmcTr = $gcmcBox.TextFrame.TextRange
$gcmcTr.Paragraphs(5).Runs(1).Text = "Use a encoder network to calculate node embeddings with Graph convolution network."
$gcmcTr.Paragraphs(6).Runs(1).Text = "Use a decoder network to predict link label with node embeddings."

# ---------------------------------------------------------------------
# 3) Slide 4 - relabel the user ("u0") and item ("i0") node shapes.
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)

$userLabels = @{
    2 = "u1"
    3 = "u7"
    4 = "u6"
    5 = "u5"
    6 = "u4"
    7 = "u3"
    8 = "u2"
}
foreach ($idx in $userLabels.Keys) {
    $slide4.Shapes.Item($idx).TextFrame.TextRange.Text = $userLabels[$idx]
}

$itemLabels = @{
    10 = "i1"
    11 = "i7"
    12 = "i6"
    13 = "i5"
    14 = "i4"
    15 = "i3"
    16 = "i2"
}
foreach ($idx in $itemLabels.Keys) {
    $slide4.Shapes.Item($idx).TextFrame.TextRange.Text = $itemLabels[$idx]
}
